$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Description text shared by the three new Deal-Chat test cases (row 60..62, column B)
$chatDescription = @"
Verify the Chat feature in Deals. 1) Enter valid shipper user id and Password and click Login button.
2) Click on Deals menu.
3) Select any deal in opportunity tab.
4) Click on the deal cards in the Opportunity tab.
5) Click on the Carrier contact name in detailed view.
6) Enter message then  Click on Send icon.
7) Enter valid carrier user id and Password and click Login button in another  Carrier User.
8) Click on Chat Icon.
9) Ensure message is displayed notification in Chat window.
"@

$chatResult = "Chat verified Successfully"

# Insert three new rows (60, 61, 62) below the existing data, copying row 37's
# formatting (A/C/D vertically-centered, B vertically-centered + wrap, 165pt tall)
# so the new rows pick up the same cell styles already used elsewhere in the sheet.
$ws.Rows.Item(37).Copy()
$ws.Rows.Item(60).Insert(-4121)
$ws.Rows.Item(37).Copy()
$ws.Rows.Item(61).Insert(-4121)
$ws.Rows.Item(37).Copy()
$ws.Rows.Item(62).Insert(-4121)

# Row 60: Deals_Chat_ShipperUser_TC001
$ws.Range("A60").Value = "Deals_Chat_ShipperUser_TC001"
$ws.Range("B60").Value = $chatDescription
$ws.Range("C60").Value = "Yes"
$ws.Range("D60").Value = $chatResult

# Row 61: Deals_Chat_ShipperAdmin_TC002
$ws.Range("A61").Value = "Deals_Chat_ShipperAdmin_TC002"
$ws.Range("B61").Value = $chatDescription
$ws.Range("C61").Value = "Yes"
$ws.Range("D61").Value = $chatResult

# Row 62: Deals_Chat_CarrierUser_TC003
$ws.Range("A62").Value = "Deals_Chat_CarrierUser_TC003"
$ws.Range("B62").Value = $chatDescription
$ws.Range("C62").Value = "Yes"
$ws.Range("D62").Value = $chatResult

# Column B needs vertical-center + wrap (matches style used elsewhere for the
# description column); the copied row only had wrap, so fix vertical alignment.
$ws.Range("B60").VerticalAlignment = -4108
$ws.Range("B61").VerticalAlignment = -4108
$ws.Range("B62").VerticalAlignment = -4108

# Match row height of the other long-description rows.
$ws.Rows.Item(60).RowHeight = 165
$ws.Rows.Item(61).RowHeight = 165
$ws.Rows.Item(62).RowHeight = 165

# Update the view: scrolled down to the new rows, with A60:D62 selected.
$ws.Range("A60:D62").Select()
